$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (updates sheet name which is reflected in the title/tab)
$ws.Name = "Through 2022-11-04"

# Update the header label in I1 (shared string)
$ws.Range("I1").Value = "2022 (through 11-04)"

# Update November row (row 12) value for 2022 column (I)
$ws.Range("I12").Value = 13

# Update Total row (row 14) value for 2022 column (I)
$ws.Range("I14").Value = 1413
